$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the A/B/C/D merged ranges from row 16 down to row 20
# to cover the four newly-added part rows (17-20).
$ws.Range("A8:A16").UnMerge()
$ws.Range("B8:B16").UnMerge()
$ws.Range("C8:C16").UnMerge()
$ws.Range("D8:D16").UnMerge()
$ws.Range("A8:A20").Merge()
$ws.Range("B8:B20").Merge()
$ws.Range("C8:C20").Merge()
$ws.Range("D8:D20").Merge()

# Row 8
$ws.Range("C8").Value = 'MLK_PMT_10102_-_V-002'
$ws.Range("D8").Value = 'Air Receiver'
$ws.Range("E8").Value = 'PLATE 1571x600x5mm SHELL TO BE ROLLED'
$ws.Range("G8").Value = 'DMS0'
$ws.Range("H8").Value = 'Stainless Steel'
$ws.Range("I8").Value = 'SA 240 M'
$ws.Range("J8").Value = '316L'
$ws.Range("K8").Value = 'N/A'
$ws.Range("L8").Value = '100 °C'
$ws.Range("M8").Value = '1.1 BarG'
$ws.Range("N8").Value = '100 °C'
$ws.Range("O8").Value = '1.0 BarG'

# Row 9
$ws.Range("E9").Value = 'DISH HEAD 8 THK TO BE FORMED 1D 2:1 TYPE'
$ws.Range("G9").Value = 'DMS0'
$ws.Range("H9").Value = 'Stainless Steel'
$ws.Range("I9").Value = 'SA 240 M'
$ws.Range("J9").Value = '316L'
$ws.Range("K9").Value = 'N/A'
$ws.Range("L9").Value = '100 °C'
$ws.Range("M9").Value = '1.1 BarG'
$ws.Range("N9").Value = '100 °C'
$ws.Range("O9").Value = '1.0 BarG'

# Row 10
$ws.Range("E10").Value = 'SEAMLESS PIPE DN50 x 87 SCH 40S'
$ws.Range("G10").Value = 'DMS0'
$ws.Range("H10").Value = 'Stainless Steel'
$ws.Range("I10").Value = 'SA 312 M TP'
$ws.Range("J10").Value = '316L'
$ws.Range("K10").Value = 'N/A'
$ws.Range("L10").Value = '100 °C'
$ws.Range("M10").Value = '1.1 BarG'
$ws.Range("N10").Value = '100 °C'
$ws.Range("O10").Value = '1.0 BarG'

# Row 11
$ws.Range("E11").Value = 'SEAMLESS PIPE DN25 x 100 SCH 40S'
$ws.Range("G11").Value = 'DMS0'
$ws.Range("H11").Value = 'Stainless Steel'
$ws.Range("I11").Value = 'SA 312 M TP'
$ws.Range("J11").Value = '316L'
$ws.Range("K11").Value = 'N/A'
$ws.Range("L11").Value = '100 °C'
$ws.Range("M11").Value = '1.1 BarG'
$ws.Range("N11").Value = '100 °C'
$ws.Range("O11").Value = '1.0 BarG'

# Row 12
$ws.Range("E12").Value = 'SEAMLESS PIPE DN50 x 112 SCH 40S'
$ws.Range("G12").Value = 'DMS0'
$ws.Range("H12").Value = 'Stainless Steel'
$ws.Range("I12").Value = 'SA 312 M TP'
$ws.Range("J12").Value = '316L'
$ws.Range("K12").Value = 'N/A'
$ws.Range("L12").Value = '100 °C'
$ws.Range("M12").Value = '1.1 BarG'
$ws.Range("N12").Value = '100 °C'
$ws.Range("O12").Value = '1.0 BarG'

# Row 13
$ws.Range("E13").Value = 'FLANGE DN50 CLASS 150 WNRF SCH 40S'
$ws.Range("G13").Value = 'DMS0'
$ws.Range("H13").Value = 'Stainless Steel'
$ws.Range("I13").Value = 'SA 182 M'
$ws.Range("J13").Value = 'F316L'
$ws.Range("K13").Value = 'N/A'
$ws.Range("L13").Value = '100 °C'
$ws.Range("M13").Value = '1.1 BarG'
$ws.Range("N13").Value = '100 °C'
$ws.Range("O13").Value = '1.0 BarG'

# Row 14
$ws.Range("E14").Value = 'FLANGE DN25 CLASS 150 WNRF SCH 40S'
$ws.Range("G14").Value = 'DMS0'
$ws.Range("H14").Value = 'Stainless Steel'
$ws.Range("I14").Value = 'SA 182 M'
$ws.Range("J14").Value = 'F316L'
$ws.Range("K14").Value = 'N/A'
$ws.Range("L14").Value = '100 °C'
$ws.Range("M14").Value = '1.1 BarG'
$ws.Range("N14").Value = '100 °C'
$ws.Range("O14").Value = '1.0 BarG'

# Row 15
$ws.Range("E15").Value = 'FLANGE DN50 CLASS 150 WNRF SCH 40S'
$ws.Range("G15").Value = 'DMS0'
$ws.Range("H15").Value = 'Stainless Steel'
$ws.Range("I15").Value = 'SA 182 M'
$ws.Range("J15").Value = 'F316L'
$ws.Range("K15").Value = 'N/A'
$ws.Range("L15").Value = '100 °C'
$ws.Range("M15").Value = '1.1 BarG'
$ws.Range("N15").Value = '100 °C'
$ws.Range("O15").Value = '1.0 BarG'

# Row 16
$ws.Range("E16").Value = 'BASE PLATE 150 x 150 x 5 THK'
$ws.Range("G16").Value = 'DMS0'
$ws.Range("H16").Value = 'Stainless Steel'
$ws.Range("I16").Value = 'SA 240 M'
$ws.Range("J16").Value = 'Gr. 304'
$ws.Range("K16").Value = 'N/A'
$ws.Range("L16").Value = '100 °C'
$ws.Range("M16").Value = '1.1 BarG'
$ws.Range("N16").Value = '100 °C'
$ws.Range("O16").Value = '1.0 BarG'

# Row 17
$ws.Range("E17").Value = 'LIFTING LUG PLATE 190 x 80 x 6 THK'
$ws.Range("G17").Value = 'DMS0'
$ws.Range("H17").Value = 'Stainless Steel'
$ws.Range("I17").Value = 'SA 240 M'
$ws.Range("J17").Value = '316L'
$ws.Range("K17").Value = 'N/A'
$ws.Range("L17").Value = '100 °C'
$ws.Range("M17").Value = '1.1 BarG'
$ws.Range("N17").Value = '100 °C'
$ws.Range("O17").Value = '1.0 BarG'

# Row 18
$ws.Range("E18").Value = 'DOUBLER PLATE 120 x 150 x 12.7 THK'
$ws.Range("G18").Value = 'DMS0'
$ws.Range("H18").Value = 'Stainless Steel'
$ws.Range("I18").Value = 'SA 240 M'
$ws.Range("J18").Value = 'Gr. 304'
$ws.Range("K18").Value = 'N/A'
$ws.Range("L18").Value = '100 °C'
$ws.Range("M18").Value = '1.1 BarG'
$ws.Range("N18").Value = '100 °C'
$ws.Range("O18").Value = '1.0 BarG'

# Row 19
$ws.Range("E19").Value = 'EQUAL ANGLE BAR 3" x 3" x 1/4" THK'
$ws.Range("G19").Value = 'DMS0'
$ws.Range("H19").Value = 'Stainless Steel'
$ws.Range("I19").Value = 'SA 240 M'
$ws.Range("J19").Value = 'Gr. 304'
$ws.Range("K19").Value = 'N/A'
$ws.Range("L19").Value = '100 °C'
$ws.Range("M19").Value = '1.1 BarG'
$ws.Range("N19").Value = '100 °C'
$ws.Range("O19").Value = '1.0 BarG'

# Row 20
$ws.Range("E20").Value = 'EARTHING LUG'
$ws.Range("G20").Value = 'DMS0'
$ws.Range("H20").Value = 'Stainless Steel'
$ws.Range("I20").Value = 'SA 240 M'
$ws.Range("J20").Value = 'Gr. 304'
$ws.Range("K20").Value = 'N/A'
$ws.Range("L20").Value = '100 °C'
$ws.Range("M20").Value = '1.1 BarG'
$ws.Range("N20").Value = '100 °C'
$ws.Range("O20").Value = '1.0 BarG'
